# Applies the "Updated cryptos list" data refresh to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / label / link updates (safe as literal strings) ---
$ws.Range('D2').Value = '68.394.34'
$ws.Range('E2').Value = '  +1.61%  '
$ws.Range('D3').Value = '3.932.12'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('E5').Value = '  +3.10%  '
$ws.Range('E6').Value = '  +1.19%  '
$ws.Range('E7').Value = '  +0.64%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +0.23%  '
$ws.Range('E10').Value = '  +2.72%  '
$ws.Range('E11').Value = '  +4.21%  '
$ws.Range('E12').Value = '  -1.06%  '
$ws.Range('E13').Value = '  +3.39%  '
$ws.Range('D14').Value = '4.560.14'
$ws.Range('E14').Value = '  -0.27%  '
$ws.Range('B15').Value = 'Uniswap'
$ws.Range('C15').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('E15').Value = '  -4.32%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.915.59'
$ws.Range('E16').Value = '  -0.80%  '
$ws.Range('E17').Value = '  -0.67%  '
$ws.Range('E18').Value = '  +0.80%  '
$ws.Range('E19').Value = '  -1.86%  '
$ws.Range('D20').Value = '68.490.52'
$ws.Range('E20').Value = '  +1.40%  '
$ws.Range('E21').Value = '  +1.93%  '
$ws.Range('E22').Value = '  +4.10%  '
$ws.Range('E23').Value = '  +3.88%  '
$ws.Range('E24').Value = '  +1.26%  '
$ws.Range('E25').Value = '  +17.10%  '
$ws.Range('E26').Value = '  +12.35%  '
$ws.Range('E27').Value = '  +1.21%  '
$ws.Range('E28').Value = '  +1.08%  '
$ws.Range('E29').Value = '  +1.30%  '
$ws.Range('E30').Value = '  -1.38%  '
$ws.Range('E31').Value = '  +1.14%  '
$ws.Range('E32').Value = '  -1.51%  '
$ws.Range('E33').Value = '  +3.78%  '
$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('E34').Value = '  -1.68%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('E35').Value = '  +15.25%  '
$ws.Range('D36').Value = '0.0₃0883'
$ws.Range('E36').Value = '  +13.08%  '
$ws.Range('E37').Value = '  +5.88%  '
$ws.Range('E38').Value = '  -1.66%  '
$ws.Range('E39').Value = '  +18.07%  '
$ws.Range('B40').Value = 'Dai'
$ws.Range('C40').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('E41').Value = '  +15.92%  '
$ws.Range('E42').Value = '  +6.38%  '
$ws.Range('E43').Value = '  +1.05%  '
$ws.Range('E44').Value = '  +4.76%  '
$ws.Range('E45').Value = '  +0.85%  '
$ws.Range('E46').Value = '  +0.07%  '
$ws.Range('D47').Value = '0.0₆0355'
$ws.Range('E47').Value = '  +40.10%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('E48').Value = '  +2.47%  '
$ws.Range('B49').Value = 'LidoDAOToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('E49').Value = '  -0.94%  '
$ws.Range('E50').Value = '  -1.92%  '
$ws.Range('E51').Value = '  -0.81%  '

# --- Price column updates that look like plain numbers: force text so
#     Excel does not auto-convert them to numeric values, then restore
#     the original (unstyled) cell formatting. ---
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '486.43'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.78'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.627'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.734'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.00'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.71'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.52'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.99'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '442.33'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.52'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '15.12'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.43'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.28'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.50'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.63'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '39.06'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.84'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '717.61'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.75'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.130'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '42.25'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.21'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '61.34'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.397'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.00'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.23'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.95'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.142'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.26'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.41'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.15'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '145.89'

foreach ($cellRef in @('D5', 'D6', 'D7', 'D9', 'D12', 'D13', 'D15', 'D18', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D34', 'D35', 'D37', 'D39', 'D40', 'D41', 'D42', 'D44', 'D45', 'D48', 'D49', 'D50', 'D51')) {
    $ws.Range($cellRef).ClearFormats()
}
